$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Apartment at Sayat-Nova Street"
$ws.Range("C10").Value = "Apartment"
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 40.169242536025401
$ws.Range("F10").Value = 44.5038657799916
$ws.Range("G10").Value = "Kentron"

$ws.Range("B17").Select()
